$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (Pin 19): "BLANK" -> "GND (Jumper)"
$ws.Range("B21").Value2 = "GND (Jumper)"

# Row 28 (Pin 26): "BLANK" -> "5V"
$ws.Range("B28").Value2 = "5V"

# Row 17 (Pin 15): "BLANK" -> "One Wire Bread Board", add comment "Thermistor in the relay"
$ws.Range("B17").Value2 = "One Wire Bread Board"
$ws.Range("C17").Value2 = "Thermistor in the relay"
$ws.Rows.Item(17).RowHeight = 30

# Update the active selection to match where the author ended up (D16)
[void]$ws.Range("D16").Select()
